$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header labels to the version-specific labels.
$fields = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fields.Length; $i++) {
    $oldCol = $i + 1          # A..J
    $newCol = $i + 12         # L..U
    $ws.Cells.Item(1, $oldCol).Value = "$($fields[$i])_FV2404"
    $ws.Cells.Item(1, $newCol).Value = "$($fields[$i])_FV2410"
}

# Turn the used range into an Excel table ("Table1") with an autofilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row and select the first cell below it.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
